$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (refreshed by GitHub Actions bot),
# including re-ranked rows 30-31 (Hedera / InternetComputer swap) and 46-48
# (Aave / Mantle / BabyDogeCoin rotation). Column D ("Price") values are forced to
# text format first since some of them (e.g. "1.005") would otherwise be
# auto-converted to numbers by Excel, which does not match the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.064.11"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.57"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.10"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6195"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07418"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2913"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.93"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07700"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.832.48"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.970"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6687"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.35"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009100"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.868"
$ws.Range("E17").Value = "  -2.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.096.84"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.083.74"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "235.29"
$ws.Range("E20").Value = "  +5.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.56"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.148"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.008"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.24"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1415"
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.493"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.76"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.493"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.110"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05549"
$ws.Range("E31").Value = "  -6.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.115"
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.840"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7403"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.136"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.654"
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.814"
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.203.35"
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.408"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9082"
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.010"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.973.12"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000123"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "64.64"
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5119"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4023"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.109"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05809"
$ws.Range("E51").Value = "  -0.34%  "
